$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.676.36"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +4.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.226.27"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.06%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.14"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.83%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.402"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.73"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.104"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.552.87"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.69"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.87"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.801"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.58"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.222.64"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.498.84"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.10"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0900"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.08"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.61"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +10.51%  "
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.52"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.34"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.141"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.99"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("D31").ClearFormats()
$ws.Range("E32").Value = "  -2.40%  "
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.97"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.64"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0627"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.60"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.14%  "
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000246"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +29.57%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0240"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.82"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.62"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +8.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0988"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +6.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "99.37"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.09%  "
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.467.56"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.54"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.77"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("E51").Value = "  -0.94%  "
